# Update countries & provincias Spain
# - Re-rank a few countries whose case counts overtook their neighbours
#   (Kuwait/Hungria, Afganistan/Hong Kong/Camerun, Etiopia/Gabon), carrying
#   each country's fresh stats along with it.
# - Refresh the "Datos actualizados" timestamp cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 11:52"

# --- Row 46: Malasia (name unchanged, stats refreshed) ---
$ws.Range("B46").Value = 5425
$ws.Range("C46").Value = 36
$ws.Range("D46").Value = 3295
$ws.Range("E46").Value = 2041
$ws.Range("F46").Value = 45

# --- Rows 61-65: Grecia / Hungria / Kuwait / Barein / Croacia ---
# Kuwait's total now exceeds Hungria's, so the two swap places (row 62 <-> 63).
$ws.Range("A62").Value = "Kuwait"
$ws.Range("B62").Value = 1995
$ws.Range("C62").Value = 80
$ws.Range("D62").Value = 367
$ws.Range("E62").Value = 1619
$ws.Range("F62").Value = 39
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 9

$ws.Range("A63").Value = "Hungria"
$ws.Range("B63").Value = 1984
$ws.Range("C63").Value = 68
$ws.Range("D63").Value = 267
$ws.Range("E63").Value = 1518
$ws.Range("F63").Value = 60
$ws.Range("G63").Value = 10
$ws.Range("H63").Value = 199

$ws.Range("B64").Value = 1895
$ws.Range("C64").Value = 14
$ws.Range("D64").Value = 769
$ws.Range("E64").Value = 1119

# --- Rows 81-85: Cuba / Hong Kong / Camerun / Afganistan / Bulgaria ---
# Afganistan jumps ahead of Hong Kong and Camerun (row 82), pushing them
# down one place each (rows 83, 84).
$ws.Range("A82").Value = "Afganistan"
$ws.Range("C82").Value = 30
$ws.Range("D82").Value = 131
$ws.Range("E82").Value = 862
$ws.Range("F82").Value = 7
$ws.Range("H82").Value = 33

$ws.Range("A83").Value = "Hong Kong"
$ws.Range("B83").Value = 1026
$ws.Range("D83").Value = 602
$ws.Range("E83").Value = 420
$ws.Range("F83").Value = 8
$ws.Range("H83").Value = 4

$ws.Range("A84").Value = "Camerun"
$ws.Range("B84").Value = 1017
$ws.Range("D84").Value = 305
$ws.Range("E84").Value = 670
$ws.Range("F84").Value = 33
$ws.Range("H84").Value = 42

# --- Row 120: Vietnam (name unchanged, stats refreshed) ---
$ws.Range("D120").Value = 207
$ws.Range("E120").Value = 61

# --- Row 133: Brunei (name unchanged, stats refreshed) ---
$ws.Range("D133").Value = 116
$ws.Range("E133").Value = 21

# --- Rows 138-142: Birmania / Gabon / Etiopia / Aruba ---
# Etiopia overtakes Gabon (row 139), pushing Gabon down to row 140.
$ws.Range("A139").Value = "Etiopia"
$ws.Range("B139").Value = 111
$ws.Range("C139").Value = 3
$ws.Range("D139").Value = 16
$ws.Range("E139").Value = 92
$ws.Range("F139").Value = 1
$ws.Range("H139").Value = 3

$ws.Range("A140").Value = "Gabon"
$ws.Range("B140").Value = 109
$ws.Range("D140").Value = 7
$ws.Range("E140").Value = 101
$ws.Range("F140").Value = 0
$ws.Range("H140").Value = 1
